# Weekly update of "Fruta / hortaliza" price records.
# Two new weekly records (2023-10-30, serial 45229) are inserted at the top of the
# historical block (rows 19-20), pushing the previously existing rows 19-25 down to
# rows 21-27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 19, shifting rows 19:25 down to 21:27.
$ws.Rows("19:20").Insert()

# ---- New row 19 : Primera calidad, new weekly record ----
$ws.Cells.Item(19, 1).Value = 10
$ws.Cells.Item(19, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(19, 3).Value = "La Araucanía"
$ws.Cells.Item(19, 4).Value = 45229
$ws.Cells.Item(19, 5).Value = 9
$ws.Cells.Item(19, 6).Value = "Fruta"
$ws.Cells.Item(19, 7).Value = 100104
$ws.Cells.Item(19, 8).Value = "Frutos de pepita"
$ws.Cells.Item(19, 9).Value = 100104004
$ws.Cells.Item(19, 10).Value = "Níspero"
$ws.Cells.Item(19, 11).Value = "Californiana(o)"
$ws.Cells.Item(19, 12).Value = "Primera"
$ws.Cells.Item(19, 13).Value = 140
$ws.Cells.Item(19, 14).Value = 32000
$ws.Cells.Item(19, 15).Value = 32000
$ws.Cells.Item(19, 16).Value = 32000
$ws.Cells.Item(19, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(19, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(19, 19).Value = 3200
$ws.Cells.Item(19, 20).Value = 10

# ---- New row 20 : Segunda calidad, new weekly record ----
$ws.Cells.Item(20, 1).Value = 10
$ws.Cells.Item(20, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(20, 3).Value = "La Araucanía"
$ws.Cells.Item(20, 4).Value = 45229
$ws.Cells.Item(20, 5).Value = 9
$ws.Cells.Item(20, 6).Value = "Fruta"
$ws.Cells.Item(20, 7).Value = 100104
$ws.Cells.Item(20, 8).Value = "Frutos de pepita"
$ws.Cells.Item(20, 9).Value = 100104004
$ws.Cells.Item(20, 10).Value = "Níspero"
$ws.Cells.Item(20, 11).Value = "Californiana(o)"
$ws.Cells.Item(20, 12).Value = "Segunda"
$ws.Cells.Item(20, 13).Value = 80
$ws.Cells.Item(20, 14).Value = 20000
$ws.Cells.Item(20, 15).Value = 20000
$ws.Cells.Item(20, 16).Value = 20000
$ws.Cells.Item(20, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(20, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(20, 19).Value = 2000
$ws.Cells.Item(20, 20).Value = 10
